$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the filled-row formatting (bold font, date number format) from row 4
# into the Name/Start/Finish columns of row 8 before setting the new values.
$ws.Range("A4:C4").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Row 8: new "Categories" task done by Nhat Kha, not yet responsive
$ws.Range("A8").Value = "Categories"
$ws.Range("B8").Value = 45063
$ws.Range("C8").Value = 45064
$ws.Range("D8").Value = 0.9
$ws.Range("E8").Value = "Nhật Kha"
$ws.Range("J8").Value = "Chưa responsive"

$ws.Range("J8").Select()

$wb.Save()
